$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(3, 6).Value = 1513
$ws1.Cells.Item(4, 6).Value = 792
$ws1.Cells.Item(5, 6).Value = 202
$ws1.Cells.Item(6, 6).Value = 54
$ws1.Cells.Item(7, 6).Value = 1097
$ws1.Cells.Item(8, 6).Value = 679
$ws1.Cells.Item(9, 6).Value = 757
$ws1.Cells.Item(10, 6).Value = 1330
$ws1.Cells.Item(12, 6).Value = 1006
$ws1.Cells.Item(13, 6).Value = 17
$ws1.Cells.Item(14, 6).Value = 53
$ws1.Cells.Item(15, 6).Value = 179
$ws1.Cells.Item(16, 6).Value = 40
$ws1.Cells.Item(17, 6).Value = 409
$ws1.Cells.Item(20, 6).Value = 527
$ws1.Cells.Item(21, 6).Value = 546
$ws1.Cells.Item(22, 6).Value = 736
$ws1.Cells.Item(23, 6).Value = 221
$ws1.Cells.Item(24, 6).Value = 159
$ws1.Cells.Item(25, 6).Value = 361

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 978
$ws2.Cells.Item(5, 6).Value = 233
$ws2.Cells.Item(8, 6).Value = 59
$ws2.Cells.Item(9, 6).Value = 579
$ws2.Cells.Item(10, 6).Value = 45
$ws2.Cells.Item(11, 6).Value = 12

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 175

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 175
$ws4.Cells.Item(4, 6).Value = 1513
$ws4.Cells.Item(6, 6).Value = 792
$ws4.Cells.Item(7, 6).Value = 202
$ws4.Cells.Item(8, 6).Value = 978
$ws4.Cells.Item(9, 6).Value = 54
$ws4.Cells.Item(10, 6).Value = 1097
$ws4.Cells.Item(11, 6).Value = 679
$ws4.Cells.Item(12, 6).Value = 757
$ws4.Cells.Item(13, 6).Value = 1330
$ws4.Cells.Item(15, 6).Value = 1006
$ws4.Cells.Item(16, 6).Value = 17
$ws4.Cells.Item(17, 6).Value = 53
$ws4.Cells.Item(18, 6).Value = 179
$ws4.Cells.Item(19, 6).Value = 40
$ws4.Cells.Item(20, 6).Value = 409
$ws4.Cells.Item(22, 6).Value = 233
$ws4.Cells.Item(28, 6).Value = 527
$ws4.Cells.Item(29, 6).Value = 546
$ws4.Cells.Item(30, 6).Value = 736
$ws4.Cells.Item(31, 6).Value = 221
$ws4.Cells.Item(32, 6).Value = 59
$ws4.Cells.Item(33, 6).Value = 159
$ws4.Cells.Item(34, 6).Value = 579
$ws4.Cells.Item(35, 6).Value = 45
$ws4.Cells.Item(36, 6).Value = 45
$ws4.Cells.Item(37, 6).Value = 12
$ws4.Cells.Item(38, 6).Value = 361
